$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: update the "Förändrad" (Changed) date column C for all data rows (2-79) from 46066 to 46070
$ws.Range("C2:C79").Value = 46070

# Step 2: apply the row-level data corrections (Beteckning/Datum/Markägare/Area) identified by the diff
# Row 14: 'A 15832-2021' -> 'A 72700-2021'
$ws.Range("A14").Value = "A 72700-2021"
$ws.Range("B14").Value = 44545
$ws.Range("G14").Value = 2.1

# Row 15: 'A 72700-2021' -> 'A 15832-2021'
$ws.Range("A15").Value = "A 15832-2021"
$ws.Range("B15").Value = 44286
$ws.Range("G15").Value = 1.8

# Row 20: 'A 41446-2023' -> 'A 61262-2024'
$ws.Range("A20").Value = "A 61262-2024"
$ws.Range("B20").Value = 45645
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = 2.4

# Row 21: 'A 30621-2023' -> 'A 61264-2024'
$ws.Range("A21").Value = "A 61264-2024"
$ws.Range("B21").Value = 45645
$ws.Range("G21").Value = 0.9

# Row 22: 'A 30623-2023' -> 'A 41446-2023'
$ws.Range("A22").Value = "A 41446-2023"
$ws.Range("B22").Value = 45175
$ws.Range("F22").Value = "Allmännings- och besparingsskogar"
$ws.Range("G22").Value = 8.4

# Row 23: 'A 61262-2024' -> 'A 36638-2021'
$ws.Range("A23").Value = "A 36638-2021"
$ws.Range("B23").Value = 44392
$ws.Range("G23").Value = 1.8

# Row 24: 'A 61264-2024' -> 'A 51096-2024'
$ws.Range("A24").Value = "A 51096-2024"
$ws.Range("B24").Value = 45603.47366898148
$ws.Range("G24").Value = 1.3

# Row 25: 'A 51096-2024' -> 'A 20197-2023'
$ws.Range("A25").Value = "A 20197-2023"
$ws.Range("B25").Value = 45055
$ws.Range("G25").Value = 3.5

# Row 26: 'A 36638-2021' -> 'A 30621-2023'
$ws.Range("A26").Value = "A 30621-2023"
$ws.Range("B26").Value = 45112.42612268519
$ws.Range("G26").Value = 3.2

# Row 27: 'A 20197-2023' -> 'A 30623-2023'
$ws.Range("A27").Value = "A 30623-2023"
$ws.Range("B27").Value = 45112.4308912037
$ws.Range("G27").Value = 2.6

# Row 29: 'A 44720-2021' -> 'A 59615-2021'
$ws.Range("A29").Value = "A 59615-2021"
$ws.Range("B29").Value = 44491.70662037037
$ws.Range("G29").Value = 6.8

# Row 30: 'A 19420-2024' -> 'A 11771-2025'
$ws.Range("A30").Value = "A 11771-2025"
$ws.Range("B30").Value = 45727
$ws.Range("G30").Value = 0.9

# Row 31: 'A 60754-2024' -> 'A 19420-2024'
$ws.Range("A31").Value = "A 19420-2024"
$ws.Range("B31").Value = 45429.42680555556
$ws.Range("G31").Value = 3.8

# Row 32: 'A 11771-2025' -> 'A 60754-2024'
$ws.Range("A32").Value = "A 60754-2024"
$ws.Range("B32").Value = 45644
$ws.Range("G32").Value = 2.4

# Row 33: 'A 5952-2024' -> 'A 44720-2021'
$ws.Range("A33").Value = "A 44720-2021"
$ws.Range("B33").Value = 44438
$ws.Range("G33").Value = 11.4

# Row 35: 'A 59615-2021' -> 'A 5952-2024'
$ws.Range("A35").Value = "A 5952-2024"
$ws.Range("B35").Value = 45335
$ws.Range("G35").Value = 2.1

# Row 52: 'A 6734-2022' -> 'A 2850-2026'
$ws.Range("A52").Value = "A 2850-2026"
$ws.Range("B52").Value = 46038.47877314815
$ws.Range("G52").Value = 3.1

# Row 53: 'A 30634-2023' -> 'A 2800-2026'
$ws.Range("A53").Value = "A 2800-2026"
$ws.Range("B53").Value = 46038
$ws.Range("G53").Value = 2.1

# Row 54: 'A 52795-2023' -> 'A 6734-2022'
$ws.Range("A54").Value = "A 6734-2022"
$ws.Range("B54").Value = 44602.46284722222
$ws.Range("G54").Value = 32.3

# Row 55: 'A 16735-2023' -> 'A 61260-2025'
$ws.Range("A55").Value = "A 61260-2025"
$ws.Range("B55").Value = 46000
$ws.Range("G55").Value = 4.4

# Row 56: 'A 19082-2023' -> 'A 30634-2023'
$ws.Range("A56").Value = "A 30634-2023"
$ws.Range("B56").Value = 45112.45543981482
$ws.Range("F56").Value = ""
$ws.Range("G56").Value = 6.2

# Row 57: 'A 38938-2021' -> 'A 52795-2023'
$ws.Range("A57").Value = "A 52795-2023"
$ws.Range("B57").Value = 45226.35605324074
$ws.Range("F57").Value = ""
$ws.Range("G57").Value = 3.7

# Row 58: 'A 2850-2026' -> 'A 16735-2023'
$ws.Range("A58").Value = "A 16735-2023"
$ws.Range("B58").Value = 45030.71976851852
$ws.Range("G58").Value = 8.699999999999999

# Row 59: 'A 2800-2026' -> 'A 19082-2023'
$ws.Range("A59").Value = "A 19082-2023"
$ws.Range("B59").Value = 45048
$ws.Range("F59").Value = "Allmännings- och besparingsskogar"
$ws.Range("G59").Value = 6.5

# Row 60: 'A 26029-2022' -> 'A 38938-2021'
$ws.Range("A60").Value = "A 38938-2021"
$ws.Range("B60").Value = 44411.63415509259
$ws.Range("G60").Value = 1.9

# Row 61: 'A 59612-2021' -> 'A 26029-2022'
$ws.Range("A61").Value = "A 26029-2022"
$ws.Range("B61").Value = 44734
$ws.Range("F61").Value = "Allmännings- och besparingsskogar"
$ws.Range("G61").Value = 9.699999999999999

# Row 62: 'A 61260-2025' -> 'A 59612-2021'
$ws.Range("A62").Value = "A 59612-2021"
$ws.Range("B62").Value = 44491
$ws.Range("G62").Value = 3.2
